$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that were removed from the dataset (identified by their
# "even_MAG-GUT*.fa" row-label in column A), from bottom to top so the
# remaining row numbers don't shift out from under us.
$ws.Rows(34).Delete()
$ws.Rows(21).Delete()
$ws.Rows(9).Delete()
$ws.Rows(8).Delete()
$ws.Rows(3).Delete()
